# "updated activity till excel form"
# Refresh the per-match batting stats (runs, balls, fours, sixes) for
# Eoin Morgan (c) / Kolkata Knight Riders in rows 2-8.
# Values are stored as text (numbers-as-text), matching the original
# sheet, so we force Text number formatting before writing the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:F8")
$rng.NumberFormat = "@"

$data = @(
    @("15", "12", "2", "0"),
    @("40", "25", "5", "2"),
    @("39", "29", "2", "2"),
    @("68", "35", "5", "6"),
    @("17", "9",  "2", "1"),
    @("34", "23", "3", "1"),
    @("30", "34", "3", "1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 3 + $j
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
